$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 17 (sldId 317) - "Moegliche Platformen"
# ---------------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$titleSh17 = $s17.Shapes.Item(1)
$bodySh17 = $s17.Shapes.Item(2)

# Resize/position the body placeholder (it previously had no explicit xfrm).
$bodySh17.Left = 56.173228346456696
$bodySh17.Top = 142.71259842519686
$bodySh17.Width = 614.9508661417323
$bodySh17.Height = 216.54330708661416

$tr17 = $bodySh17.TextFrame.TextRange

# Paragraph 1: "Shelly" -> bold
$tr17.Paragraphs(1, 1).Font.Bold = $true

# Insert "Open Source" as a new paragraph right after "Shelly" (level 2 = lvl"1").
$insA = $tr17.Paragraphs(1, 1).InsertAfter("`rOpen Source")
$insA.Paragraphs(1, 1).IndentLevel = 2

# Insert the "Nachruestbare..." paragraph (3 runs) right after "Open Source".
$openSourcePara = $tr17.Paragraphs(2, 1)
$insB = $openSourcePara.InsertAfter("`rNachrüstbare Teile f. bestehende Infrastruktur (zB. Smarte-")
$insB.Paragraphs(1, 1).IndentLevel = 2
$runRelays = $insB.InsertAfter("Relays")
$runClose = $runRelays.InsertAfter(")")

# Paragraph "Tinkerforge" -> bold (preserves the existing err="1" flag on the run).
$tinkerPara = $tr17.Paragraphs(4, 1)
$tinkerPara.Font.Bold = $true

# Insert "Open Source" / "34 Sensoren verfuegbar" / "Eher modularer..." after Tinkerforge.
$insC = $tinkerPara.InsertAfter("`rOpen Source")
$insC.Paragraphs(1, 1).IndentLevel = 2

$openSourcePara2 = $tr17.Paragraphs(5, 1)
$insD = $openSourcePara2.InsertAfter("`r34 Sensoren verfügbar")
$insD.Paragraphs(1, 1).IndentLevel = 2

$sensorenPara = $tr17.Paragraphs(6, 1)
$insE = $sensorenPara.InsertAfter("`rEher modularer/atomarer (" + [char]0x201E + "kleinere Teile" + [char]0x201C + " => mehr selbst bauen)")
$insE.Paragraphs(1, 1).IndentLevel = 2

# Remove the old trailing empty paragraph (previously just an endParaRPr).
$lastPara = $tr17.Paragraphs(8, 1)
$lastPara.Delete()

# ---------------------------------------------------------------------------
# Slide 18 (sldId 318) - "Dashboard ..."
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$titleSh18 = $s18.Shapes.Item(1)
$bodySh18 = $s18.Shapes.Item(2)

$titleSh18.Left = 56.173151606299214
$titleSh18.Top = 45.74409448818898
$titleSh18.Width = 505.11141732283465
$titleSh18.Height = 96.96850393700787

$titleSh18.TextFrame.TextRange.Text = "Dashboard Ansichten - Gliederung"

$bodySh18.Left = 56.173228346456696
$bodySh18.Top = 142.71259842519686
$bodySh18.Width = 427.3303987007874
$bodySh18.Height = 216.54331208661415

$tr18 = $bodySh18.TextFrame.TextRange

# Paragraph 1: "Ansichten" -> "Moegliche Ansichten:" (keep its pPr/rPr formatting).
$tr18.Paragraphs(1, 1).Text = "Mögliche Ansichten:"

# Paragraph 2: "Sicherheit" -> bold + new description paragraph after it.
$secPara = $tr18.Paragraphs(2, 1)
$secPara.Text = "Sicherheit"
$secPara.Font.Bold = $true
$insSec = $secPara.InsertAfter("`rZeigt zB. ob Türen abgeschlossen sind oder ob jemand ums Haus schleicht")
$insSec.Paragraphs(1, 1).IndentLevel = 2

# Paragraph "Komfort" -> bold + new description paragraph after it.
$komfortPara = $tr18.Paragraphs(4, 1)
$komfortPara.Text = "Komfort"
$komfortPara.Font.Bold = $true
$insKomfort = $komfortPara.InsertAfter("`rzB. Temperaturen, Luftfeuchtigkeit, Lichter, etc.")
$insKomfort.Paragraphs(1, 1).IndentLevel = 2

# Paragraph "Outdoor" -> bold + new description paragraph after it.
$outdoorPara = $tr18.Paragraphs(6, 1)
$outdoorPara.Text = "Outdoor"
$outdoorPara.Font.Bold = $true
$insOutdoor = $outdoorPara.InsertAfter("`rzB. Außentemperatur, andere Geräte im Außenbereich, etc.")
$insOutdoor.Paragraphs(1, 1).IndentLevel = 2

# Paragraph "Indoor" -> bold + new description paragraph after it.
$indoorPara = $tr18.Paragraphs(8, 1)
$indoorPara.Text = "Indoor"
$indoorPara.Font.Bold = $true
$insIndoor = $indoorPara.InsertAfter("`rSonstige Geräte im Innenbereich")
$insIndoor.Paragraphs(1, 1).IndentLevel = 2
